$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. "1.001", "239.16") are written with a leading apostrophe to force
# text entry, then restored to the "Normal" style so no stray number format
# / quote-prefix formatting is left behind on the cell.

$ws.Range("D2").Value = '30.776.44'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '1.882.04'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = "'239.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D7").Value = "'0.4816"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").Value = "'0.2839"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.81%  '
$ws.Range("D9").Value = "'0.06535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").Value = '1.933.50'
$ws.Range("E10").Value = '  +2.25%  '
$ws.Range("D11").Value = "'0.07485"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("D12").Value = "'16.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("D13").Value = "'5.102"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("D14").Value = "'88.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").Value = "'0.6668"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("D16").Value = '30.736.58'
$ws.Range("E16").Value = '  +1.04%  '
$ws.Range("D17").Value = "'13.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("D18").Value = '2.247.91'
$ws.Range("E18").Value = '  +5.11%  '
$ws.Range("D19").Value = "'0.9998"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("D20").Value = "'0.000007622"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.75%  '
$ws.Range("D21").Value = "'232.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.11%  '
$ws.Range("D22").Value = "'5.290"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.44%  '
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("D24").Value = "'6.173"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = "'168.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.44%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = "'9.295"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("D27").Value = "'18.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.54%  '
$ws.Range("D28").Value = "'1.939"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").Value = "'1.415"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.06%  '
$ws.Range("D30").Value = "'0.09753"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.98%  '
$ws.Range("D31").Value = "'4.349"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.46%  '
$ws.Range("D32").Value = "'4.020"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("D33").Value = "'0.05087"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").Value = "'1.222"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.07%  '
$ws.Range("D35").Value = "'0.7533"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").Value = "'2.707"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("D37").Value = "'0.01871"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("D39").Value = "'2.093"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("D40").Value = "'0.9155"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").Value = "'106.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("D42").Value = "'0.4290"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.23%  '
$ws.Range("D43").Value = "'5.792"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D45").Value = "'7.354"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.99%  '
$ws.Range("D46").Value = "'64.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("D47").Value = "'0.1290"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.84%  '
$ws.Range("D48").Value = "'1.485"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.56%  '
$ws.Range("D49").Value = "'8.997"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("D50").Value = "'33.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.60%  '
$ws.Range("D51").Value = "'0.05664"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.91%  '
